# Scheduled market-data refresh: update cached Universalis price/profit
# columns (H-N) per Leve row across the job-specific Leve Profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 5145
$ws.Range("I11").Value = 5145
$ws.Range("K11").Value = 5145
$ws.Range("M11").Value = -5005

# Row 18
$ws.Range("H18").Value = 5999
$ws.Range("I18").Value = 3998.5
$ws.Range("K18").Value = 3998.5
$ws.Range("M18").Value = -3714.5

# Row 33
$ws.Range("H33").Value = 212.05556
$ws.Range("I33").Value = 212.05556
$ws.Range("K33").Value = 212.05556
$ws.Range("M33").Value = 16.94443999999999

# Row 58
$ws.Range("H58").Value = 8216.286
$ws.Range("I58").Value = 900
$ws.Range("J58").Value = 9435.666999999999
$ws.Range("K58").Value = 2700
$ws.Range("L58").Value = 28307.001
$ws.Range("M58").Value = -2550
$ws.Range("N58").Value = -28607.001

# Row 62
$ws.Range("H62").Value = 4175.375
$ws.Range("I62").Value = 3634
$ws.Range("K62").Value = 3634
$ws.Range("M62").Value = -3010

# Row 65
$ws.Range("H65").Value = 4175.375
$ws.Range("I65").Value = 3634
$ws.Range("K65").Value = 18170
$ws.Range("M65").Value = -15050

# Row 101
$ws.Range("H101").Value = 1000
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 16
$ws.Range("H16").Value = 1122.5
$ws.Range("I16").Value = 1122.5
$ws.Range("K16").Value = 1122.5
$ws.Range("M16").Value = -835.5

# Row 19
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

# Row 74
$ws.Range("H74").Value = 628.375
$ws.Range("I74").Value = 689.7143
$ws.Range("J74").Value = 199
$ws.Range("K74").Value = 689.7143
$ws.Range("L74").Value = 199
$ws.Range("M74").Value = 184.2857
$ws.Range("N74").Value = -1947

# Row 77
$ws.Range("H77").Value = 628.375
$ws.Range("I77").Value = 689.7143
$ws.Range("J77").Value = 199
$ws.Range("K77").Value = 3448.5715
$ws.Range("L77").Value = 995
$ws.Range("M77").Value = 919.4285
$ws.Range("N77").Value = -9731

# Row 132
$ws.Range("H132").Value = 3396.8333
$ws.Range("I132").Value = 2823.8572
$ws.Range("J132").Value = 4199
$ws.Range("K132").Value = 8471.571599999999
$ws.Range("L132").Value = 12597
$ws.Range("M132").Value = -5941.571599999999
$ws.Range("N132").Value = -17657

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1519.2
$ws.Range("I86").Value = 1336.5
$ws.Range("K86").Value = 1336.5
$ws.Range("M86").Value = -213.5

# Row 89
$ws.Range("H89").Value = 1519.2
$ws.Range("I89").Value = 1336.5
$ws.Range("K89").Value = 6682.5
$ws.Range("M89").Value = -1066.5

# Row 94
$ws.Range("H94").Value = 1015
$ws.Range("I94").Value = 876.25
$ws.Range("J94").Value = 1200
$ws.Range("K94").Value = 876.25
$ws.Range("L94").Value = 1200
$ws.Range("M94").Value = -425.25
$ws.Range("N94").Value = -2102

# Row 105
$ws.Range("H105").Value = 3000
$ws.Range("I105").Value = 3000
$ws.Range("K105").Value = 3000
$ws.Range("M105").Value = -1253

# Row 140
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 283.33334
$ws.Range("J7").Value = 649.5
$ws.Range("L7").Value = 649.5
$ws.Range("N7").Value = -875.5

# Row 31
$ws.Range("H31").Value = 4512.5
$ws.Range("I31").Value = 2833.65
$ws.Range("K31").Value = 2833.65
$ws.Range("M31").Value = -2538.65

# Row 34
$ws.Range("H34").Value = 4512.5
$ws.Range("I34").Value = 2833.65
$ws.Range("K34").Value = 2833.65
$ws.Range("M34").Value = -2631.65

# Row 107
$ws.Range("H107").Value = 1120.3
$ws.Range("I107").Value = 794.875
$ws.Range("K107").Value = 794.875
$ws.Range("M107").Value = 1125.125

# Row 127
$ws.Range("H127").Value = 66666
$ws.Range("J127").Value = 66666
$ws.Range("L127").Value = 66666
$ws.Range("N127").Value = -76586

# Row 134
$ws.Range("H134").Value = 7586.5
$ws.Range("J134").Value = 3999
$ws.Range("L134").Value = 11997
$ws.Range("N134").Value = -17067

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 28
$ws.Range("I2").Value = 17.714285
$ws.Range("K2").Value = 17.714285
$ws.Range("M2").Value = 95.285715

# Row 3
$ws.Range("H3").Value = 5500
$ws.Range("I3").Value = 3000
$ws.Range("K3").Value = 3000
$ws.Range("M3").Value = -2884

# Row 132
$ws.Range("H132").Value = 5303.5
$ws.Range("I132").Value = 5071.778
$ws.Range("J132").Value = 5998.6665
$ws.Range("K132").Value = 15215.334
$ws.Range("L132").Value = 17995.9995
$ws.Range("M132").Value = -12685.334
$ws.Range("N132").Value = -23055.9995

# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 42
$ws.Range("H42").Value = 30000
$ws.Range("J42").Value = 30000
$ws.Range("L42").Value = 30000
$ws.Range("N42").Value = -31126

# Row 43
$ws.Range("H43").Value = 15665
$ws.Range("J43").Value = 15665
$ws.Range("L43").Value = 15665
$ws.Range("N43").Value = -16051

# Row 46
$ws.Range("H46").Value = 4000
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

# Row 49
$ws.Range("H49").Value = 30000
$ws.Range("J49").Value = 30000
$ws.Range("L49").Value = 30000
$ws.Range("N49").Value = -30294

# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# Row 136
$ws.Range("H136").Value = 29748.723
$ws.Range("I136").Value = 4872.5
$ws.Range("J136").Value = 36856.215
$ws.Range("K136").Value = 14617.5
$ws.Range("L136").Value = 110568.645
$ws.Range("M136").Value = -12067.5
$ws.Range("N136").Value = -115668.645

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 4210.2856
$ws.Range("I81").Value = 4210.2856
$ws.Range("K81").Value = 8420.5712
$ws.Range("M81").Value = -7359.5712

# Row 84
$ws.Range("H84").Value = 4210.2856
$ws.Range("I84").Value = 4210.2856
$ws.Range("K84").Value = 42102.856
$ws.Range("M84").Value = -36798.856

# Row 113
$ws.Range("H113").Value = 709.8
$ws.Range("J113").Value = 643.6667
$ws.Range("L113").Value = 1931.0001
$ws.Range("N113").Value = -6271.0001

# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 132
$ws.Range("H132").Value = 3566.85
$ws.Range("J132").Value = 4697.8
$ws.Range("L132").Value = 14093.4
$ws.Range("N132").Value = -19153.4

# Row 136
$ws.Range("H136").Value = 14000.5
$ws.Range("I136").Value = 14392.637
$ws.Range("K136").Value = 43177.911
$ws.Range("M136").Value = -40627.911
